$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stash existing style swatches into a scratch column (AA) before clearing the region ---
# so new cells can reuse the exact same style indices as their role-equivalent originals.
$ws.Range("A1").Copy()
$ws.Range("AA1").PasteSpecial(-4122)  # stash style "s1" from A1
$ws.Range("A2").Copy()
$ws.Range("AA2").PasteSpecial(-4122)  # stash style "s2" from A2
$ws.Range("B2").Copy()
$ws.Range("AA3").PasteSpecial(-4122)  # stash style "s3" from B2
$ws.Range("B3").Copy()
$ws.Range("AA4").PasteSpecial(-4122)  # stash style "s4" from B3
$ws.Range("B5").Copy()
$ws.Range("AA5").PasteSpecial(-4122)  # stash style "s5" from B5
$ws.Range("B13").Copy()
$ws.Range("AA6").PasteSpecial(-4122)  # stash style "s6" from B13
$ws.Range("A3").Copy()
$ws.Range("AA7").PasteSpecial(-4122)  # stash style "none" from A3
$excel.CutCopyMode = 0

# --- Clear the data region that is being rewritten ---
$ws.Range("A3:J28").Clear()

# --- Write every target row/cell, pasting the matching stashed style first, then the value ---
# Row 3
$ws.Range("AA7").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.30.0.6"
$ws.Range("AA4").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").Value = 1
$ws.Range("AA4").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = 34
$ws.Range("AA4").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = 94.90000000000001

# Row 4
$ws.Range("AA7").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.90.0.2"
$ws.Range("AA4").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = 1
$ws.Range("AA4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value = 26
$ws.Range("AA4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = 98.40000000000001

# Row 5
$ws.Range("AA7").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "Fi - 16.0 (1657)"
$ws.Range("AA4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = 1
$ws.Range("AA4").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = 105
$ws.Range("AA4").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D5").Value = 98.40000000000001

# Row 6
$ws.Range("AA1").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "Totals:"
$ws.Range("AA5").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B6").Value = 3
$ws.Range("AA5").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C6").Value = 165

# Row 12
$ws.Range("AA1").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Good Drivers (Roaming > 99.8%)"

# Row 13
$ws.Range("AA2").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "Adapter-Driver"
$ws.Range("AA3").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = "Total Samples"
$ws.Range("AA2").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("AA3").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D13").Value = "Good Roaming Calculation (%)"
$ws.Range("AA3").Copy()
$ws.Range("E13").PasteSpecial(-4122)
$ws.Range("E13").Value = "Driver Vintage"

# Row 14
$ws.Range("AA7").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("AA6").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B14").Value = 11128
$ws.Range("AA7").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("AA4").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Value = 100
$ws.Range("AA4").Copy()
$ws.Range("E14").PasteSpecial(-4122)

# Row 15
$ws.Range("AA7").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("AA6").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = 486214
$ws.Range("AA7").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("AA4").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 100
$ws.Range("AA4").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = "'2024-11-10"

# Row 16
$ws.Range("AA7").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("AA6").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Value = 11140
$ws.Range("AA7").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("AA4").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = 100
$ws.Range("AA4").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = "'2022-08-29"

# Row 17
$ws.Range("AA7").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("AA6").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Value = 14487
$ws.Range("AA7").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("AA4").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = 100
$ws.Range("AA4").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value = "'2022-05-23"

# Row 18
$ws.Range("AA7").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("AA6").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = 265400
$ws.Range("AA7").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("AA4").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 99.90000000000001
$ws.Range("AA4").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = "'2022-05-01"

# Row 19
$ws.Range("AA7").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("AA6").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B19").Value = 79953
$ws.Range("AA7").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("AA4").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D19").Value = 99.90000000000001
$ws.Range("AA4").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E19").Value = "'2021-08-18"

# Row 20
$ws.Range("AA7").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("AA6").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B20").Value = 35355
$ws.Range("AA7").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("AA4").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = 100
$ws.Range("AA4").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = "'2021-04-27"

# Row 21
$ws.Range("AA7").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("AA6").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B21").Value = 65425
$ws.Range("AA7").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("AA4").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("D21").Value = 100
$ws.Range("AA4").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("E21").Value = "'2020-08-05"

# Row 22
$ws.Range("AA7").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("AA6").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B22").Value = 117653
$ws.Range("AA7").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("AA4").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 100
$ws.Range("AA4").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'2020-01-06"

# Row 23
$ws.Range("AA7").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("AA6").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Value = 56018
$ws.Range("AA7").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("AA4").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = 100
$ws.Range("AA4").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").Value = "'2019-12-14"

$excel.CutCopyMode = 0
# --- Clean up the scratch column used for style stashing ---
$ws.Range("AA1:AA10").Clear()
